# Update source_version values:
#   Disease Ontology (row 3): v2023-05-31 -> v2023-07-20
#   Experimental Factor Ontology (row 4): v3.55.1 -> v3.56.0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

$ws.Range("E3").Value = "v2023-07-20"
$ws.Range("E4").Value = "v3.56.0"

# Update the view: scroll/selection moved from topLeftCell B1 / selection F6
# to topLeftCell C1 / selection E4.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E4").Select()
